$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.708.34"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.466.41"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'316.49"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'92.78"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").Value = "'32.69"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  +7.21%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "2.848.01"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "'15.74"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "2.522.63"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "'0.778"
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "41.684.85"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "'11.57"
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").Value = "'71.05"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "'239.68"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'2.72"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "'1.92"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'24.71"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").Value = "'9.79"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").Value = "'35.74"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("D31").Value = "'156.21"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'5.50"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "'2.50"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").Value = "'17.59"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.104"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.115"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "1.975.74"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.04"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0284"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "'9.04"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").Value = "2.701.42"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'97.08"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "'67.05"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'72.81"
$ws.Range("E51").Value = "  -0.46%  "
